$d = $word.ActiveDocument

# "Versi" + "on"  ->  single run "Version" (keeps the spellStart/spellEnd
# proofErr markers that bracket it untouched, since the Find range sits
# fully inside them).
$d.Content.Find.Execute("Version", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Version", 2)

# " 2"  ->  " 1."  (leaves the old trailing "." run, and the _GoBack
# bookmark in between, alone).
$d.Content.Find.Execute(" 2", $true, $false, $false, $false, $false, `
    $true, 1, $false, " 1.", 2)

# Drop the now-redundant "." run that trails the _GoBack bookmark; locate
# it via the bookmark itself so this isn't dependent on hard-coded offsets.
$bm = $d.Bookmarks("_GoBack")
$trailing = $d.Range($bm.End, $bm.End + 1)
if ($trailing.Text -eq ".") {
    $trailing.Delete()
}
